$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 1
    $ws.Range("F5").Value = 59
    $ws.Range("F6").Value = 1
}
